$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("25-10-2021", "26-10-2021", "27-10-2021", "28-10-2021", "29-10-2021")
$data = @(
    @(5.64, 4.85, 3.73, 3.43),
    @(5.82, 4.92, 3.79, 3.49),
    @(5.89, 4.93, 3.80, 3.52),
    @(5.91, 5.01, 3.87, 3.56),
    @(6.25, 5.10, 3.96, 3.64)
)

$startRow = 208
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
    $ws.Cells.Item($r, 5).Value = $data[$i][3]
}
